$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. ". Propust " -> " i njemu povezanim propustima. Propust "
#    (leading space stays italic like the preceding text, the rest is
#    non-italic like the original ". Propust " run)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.ClearFormatting()
$found1 = $r1.Find.Execute(". Propust ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found1) {
    $start1 = $r1.Start

    # Insert the non-italic remainder right after the leading "." boundary so
    # it naturally merges with the (already non-italic) following text
    # instead of picking up the italic formatting from the left.
    $afterDot = $d.Range($start1 + 1, $start1 + 1)
    $afterDot.InsertBefore("i njemu povezanim propustima.")

    # Remove the old leading "." and replace it with an italic space that
    # merges into the preceding italic run.
    $dotRange = $d.Range($start1, $start1 + 1)
    $dotRange.Delete()
    $spaceIns = $d.Range($start1, $start1)
    $spaceIns.InsertBefore(" ")
    $spaceRange = $d.Range($start1, $start1 + 1)
    $spaceRange.Font.Italic = 1
}

# ---------------------------------------------------------------------------
# 2. Move the _GoBack bookmark from the "...ispisuje" run into the middle of
#    "konstruiše" (between "konstrui" and "se").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$r2 = $d.Content
$r2.Find.ClearFormatting()
$found2 = $r2.Find.Execute("eksplicitno konstrui", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found2) {
    $bmPos = $r2.End
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 3. Delete the two red "Prilikom login-a / ucestala korisnicka imena"
#    paragraphs entirely.
# ---------------------------------------------------------------------------
$p1Start = -1
$p2End = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Prilikom login-a*") {
        $p1Start = $p.Range.Start
    }
    if ($t -like "*ucestala korisnicka*" -or $t -like "*u*estala korisni*ka*") {
        $p2End = $p.Range.End
    }
}
if ($p1Start -ge 0 -and $p2End -gt $p1Start) {
    $delRange = $d.Range($p1Start, $p2End)
    $delRange.Delete()
}

# ---------------------------------------------------------------------------
# 4. Replace the last (red) paragraph "je pronaci lozinku, u debug-u
#    console.log ispisuje kada sifra nije jednaka " with a new bold/italic
#    "Capture the flag #3 - CWE-209" heading followed by one empty paragraph.
# ---------------------------------------------------------------------------
$lastPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*pronaci lozinku, u debug-u console*") {
        $lastPara = $p
    }
}
if ($lastPara -ne $null) {
    $prevParaEnd = $lastPara.Range.Start

    $ins = $d.Range($prevParaEnd, $prevParaEnd)
    $ins.InsertBefore("Capture the flag #3 - CWE-209" + [char]13)

    # Re-locate the freshly inserted heading paragraph and style it.
    $headingStart = $prevParaEnd
    $headingTextEnd = $headingStart + ("Capture the flag #3 - CWE-209").Length
    $headingPara = $d.Range($headingStart, $headingTextEnd).Paragraphs.Item(1)

    $headingWhole = $headingPara.Range
    $headingWhole.Font.Bold = 1
    $headingWhole.Font.Color = -16777216

    $cweStart = $headingStart + ("Capture the flag #3 - ").Length
    $cweEnd = $headingTextEnd
    $cweRange = $d.Range($cweStart, $cweEnd)
    $cweRange.Font.Italic = 1

    # Remove the old red paragraph's visible text, leaving one empty
    # paragraph with the color stripped back to automatic.
    $oldStart = $headingTextEnd + 1
    $oldParaRangeProbe = $d.Range($oldStart, $oldStart)
    $oldPara = $oldParaRangeProbe.Paragraphs.Item(1)
    $oldPara.Range.Font.Color = -16777216

    $oldTextStart = $oldPara.Range.Start
    $oldTextEnd = $oldPara.Range.End
    $oldTextRange = $d.Range($oldTextStart, $oldTextEnd - 1)
    if ($oldTextRange.Start -lt $oldTextRange.End) {
        $oldTextRange.Delete()
    }
}
